# Update the cached "datetimeFigureOut" footer field text (date placeholder)
# on the slide master and on every slide layout from 2007-03-29 to 2007-04-01.
# This mirrors what PowerPoint does internally when it recalculates the
# auto date field and re-caches its displayed text across the deck
# (e.g. on a "Save As PowerPoint 97-2003" round trip).

$p = $ppt.ActivePresentation
$newDate = "2007-04-01"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $phType = $null
        try { $phType = $sh.PlaceholderFormat.Type } catch { $phType = $null }
        if ($phType -eq $ppPlaceholderDate) {
            if ($sh.TextFrame.TextRange.Text -ne $newDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide master footer date placeholder.
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout's footer date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $lay = $layouts.Item($li)
    Update-DatePlaceholder $lay.Shapes
}
